$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-summary"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-id}
"
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/category}
"
$elements.Range("Q9").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/category"
$elements.Range("Y11").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/insight-category-values"
$elements.Range("Q12").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-summary"
